$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 714 (pushes existing rows 714+ down to 715+)
$ws.Rows(714).Insert()

# Populate the newly inserted row with the new data point for 2026/01/24
# Leading apostrophe forces text storage so the date-shaped string is not
# reinterpreted as a date serial number (matches the other rows in column A).
$ws.Range("A714").Value = "'2026/01/24"
$ws.Range("B714").Value = "土"
$ws.Range("C714").Value = 13
$ws.Range("D714").Value = 130
